# The post stored in row 189 ("「過去を忘れろ。だが教訓は忘れるな」...") was
# removed from the source data. Delete that entire row so that every
# subsequent row shifts up by one, and the sheet dimension shrinks
# from A1:C281 to A1:C280.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(189).Delete()
